$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "-4,7"
$ws.Range("E2").Value = "13,7"
$ws.Range("D3").Value = "0,6"
$ws.Range("E3").Value = "4,7"
$ws.Range("D4").Value = "1,2"
$ws.Range("E4").Value = "5,9"
$ws.Range("D5").Value = "-0,6"
$ws.Range("E5").Value = "6,3"
$ws.Range("D6").Value = "-4,8"
$ws.Range("E6").Value = "0,3"
$ws.Range("D7").Value = "-2,1"
$ws.Range("E7").Value = "3,3"
$ws.Range("D8").Value = "3,4"
$ws.Range("E8").Value = "4,2"
$ws.Range("D9").Value = "2,3"
$ws.Range("E9").Value = "4,8"
$ws.Range("D10").Value = "4,6"
$ws.Range("E10").Value = "9,8"
$ws.Range("D11").Value = "9,8"
$ws.Range("E11").Value = "10,6"
$ws.Range("D12").Value = "10,6"
$ws.Range("E12").Value = "11,4"
$ws.Range("D13").Value = "11,3"
$ws.Range("E13").Value = "12,2"
$ws.Range("D14").Value = "11,3"
$ws.Range("E14").Value = "12,3"
$ws.Range("D15").Value = "10,2"
$ws.Range("E15").Value = "11,7"
$ws.Range("D16").Value = "9,9"
$ws.Range("E16").Value = "11,8"
$ws.Range("D17").Value = "9,5"
$ws.Range("E17").Value = "12,3"
$ws.Range("D18").Value = "8,7"
$ws.Range("E18").Value = "11,3"
$ws.Range("D19").Value = "9,7"
$ws.Range("E19").Value = "11,1"
$ws.Range("D20").Value = "8,6"
$ws.Range("E20").Value = "9,9"
$ws.Range("D21").Value = "8,7"
$ws.Range("E21").Value = "9,5"
$ws.Range("D22").Value = "9,6"
$ws.Range("E22").Value = "10,9"
$ws.Range("D23").Value = "9,6"
$ws.Range("E23").Value = "10,1"
$ws.Range("D24").Value = "9,8"
$ws.Range("E24").Value = "11,1"
$ws.Range("D25").Value = "9,5"
$ws.Range("E25").Value = "11,1"
$ws.Range("E26").Value = "11,8"
$ws.Range("D27").Value = "11,7"
$ws.Range("E27").Value = "12,9"
$ws.Range("D28").Value = "11,2"
$ws.Range("D29").Value = "9,8"
$ws.Range("E29").Value = "11,9"
$ws.Range("D30").Value = "7,7"
$ws.Range("E31").Value = "8,7"
$ws.Range("D32").Value = "8,4"
$ws.Range("E32").Value = "12,9"
$ws.Range("D33").Value = "6,9"
$ws.Range("E33").Value = "10,7"
$ws.Range("D34").Value = "7,3"
$ws.Range("D35").Value = "-5,2"
$ws.Range("E35").Value = "40,4"
$ws.Range("D36").Value = "0,7"
$ws.Range("E36").Value = "7,9"
$ws.Range("E37").Value = "2,9"
